$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Clear C1:Q1 so the sheet's used range shrinks back to A1:B1
$ws.Range("C1:Q1").Clear()

# Update the remaining two values
$ws.Range("A1").Value = 36
$ws.Range("B1").Value = 37
